$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Phase 3 RAD Non-UI Test Cases and Data: append new FEIN-mismatch rows
# covering the additional tax types for the "Existing Liability w/Notice
# Number" and "New Tax Return Amount Due" payment types.

$existingLiabilityTaxTypes = @(
    "Admissions and Amusement Tax",
    "Estate Tax",
    "Motor Fuel Tax",
    "Slots License Fee",
    "Tobacco Tax",
    "Transportation Network Services",
    "Unclaimed Property",
    "IFTA Tax"
)

$newTaxReturnTaxTypes = @(
    "Admissions and Amusement Tax",
    "Alcohol Tax",
    "Bay Restoration Fee",
    "Corporate Income Tax",
    "Estate Tax",
    "Motor Fuel Tax",
    "Sales and Use Tax",
    "Slots License Fee",
    "Tire Recycling Fee",
    "Tobacco Tax",
    "Transportation Network Services",
    "Unclaimed Property",
    "Withholding Tax"
)

$row = 14
foreach ($taxType in $existingLiabilityTaxTypes) {
    $ws1.Range("C" + $row).Value = "Y"
    $ws1.Range("D" + $row).Value = "Existing Liability w/Notice Number"
    $ws1.Range("E" + $row).Value = $taxType
    $row = $row + 1
}

foreach ($taxType in $newTaxReturnTaxTypes) {
    $ws1.Range("C" + $row).Value = "Y"
    $ws1.Range("D" + $row).Value = "New Tax Return Amount Due"
    $ws1.Range("E" + $row).Value = $taxType
    $row = $row + 1
}

# Switch active sheet/selection to reflect the final saved view state.
$ws1.Range("E27").Select()
